$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.625153034975597
$ws.Cells.Item(2, 3).Value = 0.1858221518913012
$ws.Cells.Item(2, 4).Value = 0.133060978564508
$ws.Cells.Item(2, 5).Value = 0.1137359540516316
$ws.Cells.Item(2, 6).Value = 1.471928477265259
$ws.Cells.Item(2, 9).Value = 0.8647057324363807
$ws.Cells.Item(2, 10).Value = 0.1314700272970006
$ws.Cells.Item(2, 12).Value = 0.3616644208881752
$ws.Cells.Item(2, 15).Value = 3.715882654556964

$ws.Cells.Item(3, 2).Value = 1.488498657324328
$ws.Cells.Item(3, 3).Value = 0.1663086030158638
$ws.Cells.Item(3, 4).Value = 0.1320441260592702
$ws.Cells.Item(3, 5).Value = 0.1145093657758265
$ws.Cells.Item(3, 6).Value = 1.483164092666129
$ws.Cells.Item(3, 9).Value = 0.8782526259561507
$ws.Cells.Item(3, 10).Value = 0.1329526615388834
$ws.Cells.Item(3, 12).Value = 0.3517483801026486
$ws.Cells.Item(3, 15).Value = 3.751506263169759

$ws.Cells.Item(4, 2).Value = 1.404576991406202
$ws.Cells.Item(4, 3).Value = 0.1542672123199225
$ws.Cells.Item(4, 4).Value = 0.1314520166541087
$ws.Cells.Item(4, 5).Value = 0.1150229574208987
$ws.Cells.Item(4, 6).Value = 1.490983473710273
$ws.Cells.Item(4, 9).Value = 0.8871437396526858
$ws.Cells.Item(4, 10).Value = 0.1339139388602277
$ws.Cells.Item(4, 12).Value = 0.3457626284376545
$ws.Cells.Item(4, 15).Value = 3.775907542300132

$ws.Cells.Item(5, 2).Value = 1.370376729839279
$ws.Cells.Item(5, 3).Value = 0.1493454704451835
$ws.Cells.Item(5, 4).Value = 0.1312188800539644
$ws.Cells.Item(5, 5).Value = 0.1152420000336614
$ws.Cells.Item(5, 6).Value = 1.49440134721462
$ws.Cells.Item(5, 9).Value = 0.8909109525795067
$ws.Cells.Item(5, 10).Value = 0.1343184849287083
$ws.Cells.Item(5, 12).Value = 0.3433494274880928
$ws.Cells.Item(5, 15).Value = 3.786486431882281

$ws.Cells.Item(6, 2).Value = 1.364697781428163
$ws.Cells.Item(6, 3).Value = 0.1485273355704351
$ws.Cells.Item(6, 4).Value = 0.1311806616426878
$ws.Cells.Item(6, 5).Value = 0.1152789611994214
$ws.Cells.Item(6, 6).Value = 1.494982856386819
$ws.Cells.Item(6, 9).Value = 0.8915451892554103
$ws.Cells.Item(6, 10).Value = 0.1343864339504586
$ws.Cells.Item(6, 12).Value = 0.3429502962748501
$ws.Cells.Item(6, 15).Value = 3.788281394248912

$ws.Cells.Item(7, 2).Value = 1.404115757664783
$ws.Cells.Item(7, 3).Value = 0.1542008954075698
$ws.Cells.Item(7, 4).Value = 0.1314488394215871
$ws.Cells.Item(7, 5).Value = 0.1150258720094506
$ws.Cells.Item(7, 6).Value = 1.491028631499567
$ws.Cells.Item(7, 9).Value = 0.8871939626931997
$ws.Cells.Item(7, 10).Value = 0.133919342798311
$ws.Cells.Item(7, 12).Value = 0.3457299774485136
$ws.Cells.Item(7, 15).Value = 3.776047641897307

$ws.Cells.Item(8, 2).Value = 1.578039149006145
$ws.Cells.Item(8, 3).Value = 0.1791065100528897
$ws.Cells.Item(8, 4).Value = 0.1327037025413276
$ws.Cells.Item(8, 5).Value = 0.1139946003789021
$ws.Cells.Item(8, 6).Value = 1.475611373455699
$ws.Cells.Item(8, 9).Value = 0.869257636729774
$ws.Cells.Item(8, 10).Value = 0.1319706744064177
$ws.Cells.Item(8, 12).Value = 0.358224167224634
$ws.Cells.Item(8, 15).Value = 3.727640591528427

$ws.Cells.Item(9, 2).Value = 1.918894934820344
$ws.Cells.Item(9, 3).Value = 0.2274591484837174
$ws.Cells.Item(9, 4).Value = 0.1354185305518527
$ws.Cells.Item(9, 5).Value = 0.1122787943207477
$ws.Cells.Item(9, 6).Value = 1.452688390071181
$ws.Cells.Item(9, 9).Value = 0.8386391077839228
$ws.Cells.Item(9, 10).Value = 0.1285529823018092
$ws.Cells.Item(9, 12).Value = 0.3835328042515016
$ws.Cells.Item(9, 15).Value = 3.652802128671539

$ws.Cells.Item(10, 2).Value = 2.169105577986045
$ws.Cells.Item(10, 3).Value = 0.2626752770324288
$ws.Cells.Item(10, 4).Value = 0.1375658956400656
$ws.Cells.Item(10, 5).Value = 0.1112041559391255
$ws.Cells.Item(10, 6).Value = 1.440311556566542
$ws.Cells.Item(10, 9).Value = 0.8189279029200875
$ws.Cells.Item(10, 10).Value = 0.1262873216949436
$ws.Cells.Item(10, 12).Value = 0.4026113619403446
$ws.Cells.Item(10, 15).Value = 3.610103593473127

$ws.Cells.Item(11, 2).Value = 2.282868830490315
$ws.Cells.Item(11, 3).Value = 0.2786268128067491
$ws.Cells.Item(11, 4).Value = 0.1385755415547294
$ws.Cells.Item(11, 5).Value = 0.1107554668581852
$ws.Cells.Item(11, 6).Value = 1.435652315916229
$ws.Cells.Item(11, 9).Value = 0.8105669614605269
$ws.Cells.Item(11, 10).Value = 0.1253097210514655
$ws.Cells.Item(11, 12).Value = 0.4113943097949573
$ws.Cells.Item(11, 15).Value = 3.593355146968293

$ws.Cells.Item(12, 2).Value = 2.325937589928685
$ws.Cells.Item(12, 3).Value = 0.284657144192181
$ws.Cells.Item(12, 4).Value = 0.1389625437591633
$ws.Cells.Item(12, 5).Value = 0.1105913215114143
$ws.Cells.Item(12, 6).Value = 1.434027754013144
$ws.Cells.Item(12, 9).Value = 0.8074881550827833
$ws.Cells.Item(12, 10).Value = 0.1249471488524221
$ws.Cells.Item(12, 12).Value = 0.4147349662546134
$ws.Cells.Item(12, 15).Value = 3.587398342472454

$ws.Cells.Item(13, 2).Value = 2.316662486322286
$ws.Cells.Item(13, 3).Value = 0.2833588618410374
$ws.Cells.Item(13, 4).Value = 0.1388789888739197
$ws.Cells.Item(13, 5).Value = 0.1106264170206401
$ws.Cells.Item(13, 6).Value = 1.434371412743715
$ws.Cells.Item(13, 9).Value = 0.8081473453840502
$ws.Cells.Item(13, 10).Value = 0.1250248962795524
$ws.Cells.Item(13, 12).Value = 0.4140148437110014
$ws.Cells.Item(13, 15).Value = 3.588664091021258

$ws.Cells.Item(14, 2).Value = 2.286412356000369
$ws.Cells.Item(14, 3).Value = 0.279123137450739
$ws.Cells.Item(14, 4).Value = 0.13860728709912
$ws.Cells.Item(14, 5).Value = 0.1107418470865511
$ws.Cells.Item(14, 6).Value = 1.435515859709753
$ws.Cells.Item(14, 9).Value = 0.8103119154856309
$ws.Cells.Item(14, 10).Value = 0.125279739257103
$ws.Cells.Item(14, 12).Value = 0.4116688533199095
$ws.Cells.Item(14, 15).Value = 3.592857347373382

$ws.Cells.Item(15, 2).Value = 2.267881783702592
$ws.Cells.Item(15, 3).Value = 0.2765273000628952
$ws.Cells.Item(15, 4).Value = 0.1384414689088587
$ws.Cells.Item(15, 5).Value = 0.1108133015288288
$ws.Cells.Item(15, 6).Value = 1.436235076118351
$ws.Cells.Item(15, 9).Value = 0.8116491519248896
$ws.Cells.Item(15, 10).Value = 0.1254368306267235
$ws.Cells.Item(15, 12).Value = 0.4102337803744547
$ws.Cells.Item(15, 15).Value = 3.59547605984838

$ws.Cells.Item(16, 2).Value = 2.161669487703989
$ws.Cells.Item(16, 3).Value = 0.2616314012651912
$ws.Cells.Item(16, 4).Value = 0.1375005690293278
$ws.Cells.Item(16, 5).Value = 0.111234286037444
$ws.Cells.Item(16, 6).Value = 1.440635602181942
$ws.Cells.Item(16, 9).Value = 0.8194865162553135
$ws.Cells.Item(16, 10).Value = 0.1263522773796617
$ws.Cells.Item(16, 12).Value = 0.4020394505038496
$ws.Cells.Item(16, 15).Value = 3.611252043332939

$ws.Cells.Item(17, 2).Value = 2.096494863186251
$ws.Cells.Item(17, 3).Value = 0.2524754839980403
$ws.Cells.Item(17, 4).Value = 0.1369317257235849
$ws.Cells.Item(17, 5).Value = 0.1115028255792687
$ws.Cells.Item(17, 6).Value = 1.44358400159517
$ws.Cells.Item(17, 9).Value = 0.8244497763616963
$ws.Cells.Item(17, 10).Value = 0.1269274598909287
$ws.Cells.Item(17, 12).Value = 0.3970389838449222
$ws.Cells.Item(17, 15).Value = 3.621615856110878

$ws.Cells.Item(18, 2).Value = 2.059002751585638
$ws.Cells.Item(18, 3).Value = 0.2472028182017709
$ws.Cells.Item(18, 4).Value = 0.1366076319623346
$ws.Cells.Item(18, 5).Value = 0.1116610641109048
$ws.Cells.Item(18, 6).Value = 1.445371232821373
$ws.Cells.Item(18, 9).Value = 0.827361520278906
$ws.Cells.Item(18, 10).Value = 0.1272632842407435
$ws.Cells.Item(18, 12).Value = 0.3941726508159746
$ws.Cells.Item(18, 15).Value = 3.627828616754016

$ws.Cells.Item(19, 2).Value = 2.046307716041781
$ws.Cells.Item(19, 3).Value = 0.2454164907249776
$ws.Cells.Item(19, 4).Value = 0.1364984314518196
$ws.Cells.Item(19, 5).Value = 0.1117152908591894
$ws.Cells.Item(19, 6).Value = 1.445992050116431
$ws.Cells.Item(19, 9).Value = 0.8283571721350356
$ws.Cells.Item(19, 10).Value = 0.12737784660384
$ws.Cells.Item(19, 12).Value = 0.3932038499028607
$ws.Cells.Item(19, 15).Value = 3.629975369257693

$ws.Cells.Item(20, 2).Value = 2.103433389079328
$ws.Cells.Item(20, 3).Value = 0.2534508142809671
$ws.Cells.Item(20, 4).Value = 0.1369919606042913
$ws.Cells.Item(20, 5).Value = 0.1114738477890196
$ws.Cells.Item(20, 6).Value = 1.443260679795749
$ws.Cells.Item(20, 9).Value = 0.8239155273623417
$ws.Cells.Item(20, 10).Value = 0.1268657138535596
$ws.Cells.Item(20, 12).Value = 0.3975702790152837
$ws.Cells.Item(20, 15).Value = 3.620486547815915

$ws.Cells.Item(21, 2).Value = 2.295297868254579
$ws.Cells.Item(21, 3).Value = 0.2803675509378252
$ws.Cells.Item(21, 4).Value = 0.138686966080563
$ws.Cells.Item(21, 5).Value = 0.1107077861611039
$ws.Cells.Item(21, 6).Value = 1.435175912752712
$ws.Cells.Item(21, 9).Value = 0.8096737579509217
$ws.Cells.Item(21, 10).Value = 0.1252046788323284
$ws.Cells.Item(21, 12).Value = 0.4123575290860941
$ws.Cells.Item(21, 15).Value = 3.591615218791958

$ws.Cells.Item(22, 2).Value = 2.420627649643563
$ws.Cells.Item(22, 3).Value = 0.2978997729860566
$ws.Cells.Item(22, 4).Value = 0.1398219567278005
$ws.Cells.Item(22, 5).Value = 0.1102407085991501
$ws.Cells.Item(22, 6).Value = 1.430706933339806
$ws.Cells.Item(22, 9).Value = 0.8008748814609277
$ws.Cells.Item(22, 10).Value = 0.1241635319389398
$ws.Cells.Item(22, 12).Value = 0.4221076975886291
$ws.Cells.Item(22, 15).Value = 3.574993292916758

$ws.Cells.Item(23, 2).Value = 2.353743516891598
$ws.Cells.Item(23, 3).Value = 0.2885480397073934
$ws.Cells.Item(23, 4).Value = 0.139213716324015
$ws.Cells.Item(23, 5).Value = 0.1104869276063525
$ws.Cells.Item(23, 6).Value = 1.433017498439781
$ws.Cells.Item(23, 9).Value = 0.8055243765186688
$ws.Cells.Item(23, 10).Value = 0.1247151477218618
$ws.Cells.Item(23, 12).Value = 0.4168960641976298
$ws.Cells.Item(23, 15).Value = 3.583658862588749

$ws.Cells.Item(24, 2).Value = 2.100296550960991
$ws.Cells.Item(24, 3).Value = 0.2530098948684554
$ws.Cells.Item(24, 4).Value = 0.1369647192491996
$ws.Cells.Item(24, 5).Value = 0.111486936649916
$ws.Cells.Item(24, 6).Value = 1.44340656652728
$ws.Cells.Item(24, 9).Value = 0.8241568797472887
$ws.Cells.Item(24, 10).Value = 0.1268936132122178
$ws.Cells.Item(24, 12).Value = 0.3973300539719986
$ws.Cells.Item(24, 15).Value = 3.62099631549674

$ws.Cells.Item(25, 2).Value = 1.826715976498349
$ws.Cells.Item(25, 3).Value = 0.2144318424827532
$ws.Cells.Item(25, 4).Value = 0.1346571176349585
$ws.Cells.Item(25, 5).Value = 0.1127102435519767
$ws.Cells.Item(25, 6).Value = 1.458106206644224
$ws.Cells.Item(25, 9).Value = 0.8464338436636503
$ws.Cells.Item(25, 10).Value = 0.1294344242990135
$ws.Cells.Item(25, 12).Value = 0.3766005252826687
$ws.Cells.Item(25, 15).Value = 3.67089332168166
